# Append-style refresh of the scraped job list: keep only the header row plus
# the single top-priority posting (what used to be row 6), stamp it with the
# new scrape timestamp, drop every other previously-listed job, and shrink a
# few columns that no longer need to be as wide.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Drop the old rows 3-19 entirely (their data, hyperlinks, everything).
$ws.Range("A3:H19").EntireRow.Delete()

# Row 2 becomes the posting that used to live in row 6, with a refreshed
# "fetched at" timestamp.
$ws.Range("A2").Value = "2026-02-07 06:39:11"
$ws.Range("B2").Value = "【業務自動化】国際郵便マイページの配送ラベル一括印刷の自動化ツール開発"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("G2").Value = 205
$ws.Range("H2").Value = "◆ツール,開発"

# Repoint F2's hyperlink at the new job's URL. Wiping first also clears out
# every stale hyperlink relationship left behind by the deleted rows.
$ws.Range("F2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5487449")
$ws.Range("F2").Style = "Hyperlink"

# Shrink a few columns to match the new content widths. COM's ColumnWidth
# round-trips through a pixel conversion that adds ~0.8333 characters back
# on save, so subtract it here to land exactly on the target widths.
$ws.Columns.Item(2).ColumnWidth = 37 - 5/6
$ws.Columns.Item(4).ColumnWidth = 26 - 5/6
$ws.Columns.Item(8).ColumnWidth = 12 - 5/6
